# Commit: Tue, Mar 24, 2020  5:05:51 PM
#
# The cash-flow comparison table on the "PLENARY" slide had its table
# style switched (via Table Design > Table Styles gallery) to a
# different built-in style.

$p = $ppt.ActivePresentation

# The table lives on slide 16 ("PLENARY- COMPLETE THE MISSING GAPS").
$slide = $p.Slides.Item(16)

# Find the shape that actually holds the table instead of hard-coding
# its index, so the script stays correct even if shapes get reordered.
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $tableShape = $shape
        break
    }
}

if ($tableShape -ne $null) {
    # Re-apply the table with a different style from the built-in
    # Table Styles gallery (identified by its style GUID).
    $tableShape.Table.ApplyStyle("{66E8017C-ACCE-436E-94C6-40FC3E8F40AF}")
}
